# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to match the newly scraped data (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$sheetExhibition = $wb.Worksheets.Item("展览")
$updatesExhibition = @{
    3  = 5557
    5  = 236
    9  = 4385
    11 = 818
    12 = 41
    13 = 31
    14 = 127
    15 = 144
    19 = 606
    22 = 1138
    24 = 2765
    26 = 328
}
foreach ($row in $updatesExhibition.Keys) {
    $sheetExhibition.Cells.Item($row, 6).Value = $updatesExhibition[$row]
}

# Sheet "全部类型": row -> new F value (row numbers are shifted by 1
# after row 19 relative to "展览", due to an extra row in this sheet)
$sheetAllTypes = $wb.Worksheets.Item("全部类型")
$updatesAllTypes = @{
    3  = 5557
    5  = 236
    9  = 4385
    11 = 818
    12 = 41
    13 = 31
    14 = 127
    15 = 144
    19 = 606
    23 = 1138
    25 = 2765
    27 = 328
}
foreach ($row in $updatesAllTypes.Keys) {
    $sheetAllTypes.Cells.Item($row, 6).Value = $updatesAllTypes[$row]
}
